$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Scattered single-cell value changes (imputed / removed data points) ---
$ws.Range("F2").ClearContents()            # 18.03 -> blank
$ws.Range("E6").Value = -5.7               # blank -> -5.7
$ws.Range("E8").ClearContents()            # -6.6 -> blank
$ws.Range("E18").Value = -8.5              # blank -> -8.5
$ws.Range("E20").ClearContents()           # -7.2 -> blank
$ws.Range("E23").Value = -7                # blank -> -7
$ws.Range("E25").ClearContents()           # -7.1 -> blank

# --- Remove the "RM 232" row and the "SC 92" row entirely (shrinks used range by 2 rows) ---
$ws.Rows(28).Delete()   # SC 92 row (was row 28)
$ws.Rows(26).Delete()   # RM 232 row (was row 26)

# --- Fix up the B/E/F columns for the rows that shifted up, to match the new data set ---
$ws.Range("B27").Value = -20.4             # SC 101: blank -> -20.4
$ws.Range("B28").ClearContents()           # SC 105: -19.6 -> blank
$ws.Range("B29").ClearContents()           # SC 119: -19.5 -> blank
$ws.Range("B30").Value = -19.7             # SC 120: blank -> -19.7
$ws.Range("E30").Value = -5.7              # SC 120: blank -> -5.7
$ws.Range("F30").Value = 16.89             # SC 120: blank -> 16.89
$ws.Range("B32").ClearContents()           # SC 193: -19.9 -> blank
